$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Sport"
$ws.Cells.Item(1, 2).Value = "Achievement"
$ws.Cells.Item(1, 3).Value = "Year"
$ws.Cells.Item(1, 4).Value = "Description"

$ws.Cells.Item(2, 1).Value = "Archery"
$ws.Cells.Item(2, 2).Value = "Started"
$ws.Cells.Item(2, 3).Value = 2025
$ws.Cells.Item(2, 4).Value = "I use the 35lb left-hand bow"

$ws.Cells.Item(3, 1).Value = "Basketball"
$ws.Cells.Item(3, 2).Value = "Started for Assassins"
$ws.Cells.Item(3, 3).Value = 2018
$ws.Cells.Item(3, 4).Value = "I started to learn the sport and was quite a noob"

$ws.Cells.Item(4, 1).Value = "Basketball"
$ws.Cells.Item(4, 2).Value = "Captain of the team"
$ws.Cells.Item(4, 3).Value = 2020
$ws.Cells.Item(4, 4).Value = "After many months of hard work, I was recognized for it and became a leader"

$ws.Cells.Item(5, 1).Value = "Basketball"
$ws.Cells.Item(5, 2).Value = "Runner-Up "
$ws.Cells.Item(5, 3).Value = 2020
$ws.Cells.Item(5, 4).Value = "Came 2nd Place for 3 v 3 Tournament"

$ws.Cells.Item(6, 1).Value = "Basketball"
$ws.Cells.Item(6, 2).Value = "Coach and Captain"
$ws.Cells.Item(6, 3).Value = "2023 - 2024"
$ws.Cells.Item(6, 4).Value = "Was the coach and captain of College House Residence for the interleague"

$ws.Cells.Item(7, 1).Value = "Tennis"
$ws.Cells.Item(7, 2).Value = "Runner-Up "
$ws.Cells.Item(7, 3).Value = 2019
$ws.Cells.Item(7, 4).Value = "Came 2nd Place for my division "

$ws.Cells.Item(8, 1).Value = "Table Tennis"
$ws.Cells.Item(8, 2).Value = "Hobby"
$ws.Cells.Item(8, 3).Value = "Ongoing"
$ws.Cells.Item(8, 4).Value = "Continue to grow my skills"

$ws.Cells.Item(9, 1).Value = "Squash"
$ws.Cells.Item(9, 2).Value = "Started"
$ws.Cells.Item(9, 3).Value = 2024
$ws.Cells.Item(9, 4).Value = "Tried the sport out to learn about it."

$ws.Cells.Item(10, 1).Value = "Cricket"
$ws.Cells.Item(10, 2).Value = "Started"
$ws.Cells.Item(10, 3).Value = 2019
$ws.Cells.Item(10, 4).Value = "Tried the sport out to learn about it. And nickname the Piosoner for my serves."

$ws.Cells.Item(11, 1).Value = "Baseball"
$ws.Cells.Item(11, 2).Value = "Started"
$ws.Cells.Item(11, 3).Value = 2017
$ws.Cells.Item(11, 4).Value = "Tried the sport out to learn about it."

$ws.Cells.Item(12, 1).Value = "Soccer"
$ws.Cells.Item(12, 2).Value = "Started"
$ws.Cells.Item(12, 3).Value = "2015 - 2019"
$ws.Cells.Item(12, 4).Value = "Tried the sport but I am no good."

$ws.Cells.Item(13, 1).Value = "Rugby"
$ws.Cells.Item(13, 2).Value = "Started"
$ws.Cells.Item(13, 3).Value = 2018
$ws.Cells.Item(13, 4).Value = "Not the best at it and not my favourite"

$ws.Cells.Item(14, 1).Value = "Netball"
$ws.Cells.Item(14, 2).Value = "For Fun"
$ws.Cells.Item(14, 3).Value = 2019
$ws.Cells.Item(14, 4).Value = "I played against my high school netball team and classmates, and I was a very good shooter (thanks to basketball). My basketball mates and I had to play for the court and winner gets the court."

$ws.Cells.Item(15, 1).Value = "Taekondo"
$ws.Cells.Item(15, 2).Value = "Hobby"
$ws.Cells.Item(15, 3).Value = 2024
$ws.Cells.Item(15, 4).Value = "Self-Defense and Discipline"

$ws.Cells.Item(16, 1).Value = "Tai chi"
$ws.Cells.Item(16, 2).Value = "Hobby"
$ws.Cells.Item(16, 3).Value = 2024
$ws.Cells.Item(16, 4).Value = "Discipline. To help meditate."

$ws.Cells.Item(17, 1).Value = "Boxing"
$ws.Cells.Item(17, 2).Value = "Started"
$ws.Cells.Item(17, 3).Value = 2017
$ws.Cells.Item(17, 4).Value = "Father use to teach me how to box with punching bag."

$ws.Cells.Item(18, 1).Value = "Volleyball"
$ws.Cells.Item(18, 2).Value = "Started"
$ws.Cells.Item(18, 3).Value = "2019 - 2020"
$ws.Cells.Item(18, 4).Value = "Was very good with receiver. And went to a boot camp just to get some workout in."

$ws.Cells.Item(19, 1).Value = "Badminton"
$ws.Cells.Item(19, 2).Value = "Hobby"
$ws.Cells.Item(19, 3).Value = 2024
$ws.Cells.Item(19, 4).Value = "Tried the sport out to learn about it."

$ws.Cells.Item(20, 1).Value = "Pool"
$ws.Cells.Item(20, 2).Value = "Started"
$ws.Cells.Item(20, 3).Value = 2024
$ws.Cells.Item(20, 4).Value = "I played at College House Pool table for fun."

$ws.Cells.Item(21, 1).Value = "Chess"
$ws.Cells.Item(21, 2).Value = "Started"
$ws.Cells.Item(21, 3).Value = 2016
$ws.Cells.Item(21, 4).Value = "I was taught by my father and joined a school chess team."

$ws.Cells.Item(22, 1).Value = "Chess"
$ws.Cells.Item(22, 2).Value = "Runner-Up "
$ws.Cells.Item(22, 3).Value = 2016
$ws.Cells.Item(22, 4).Value = "I won every game on my team and became a leader. But overall my team came second."

$ws.Cells.Item(23, 1).Value = "Chess"
$ws.Cells.Item(23, 2).Value = "University "
$ws.Cells.Item(23, 3).Value = "2023 - 2024"
$ws.Cells.Item(23, 4).Value = "Join College House team"

$ws.Cells.Item(24, 1).Value = "Chess"
$ws.Cells.Item(24, 2).Value = "Hobby"
$ws.Cells.Item(24, 3).Value = "Ongoing"
$ws.Cells.Item(24, 4).Value = "Average ELO Rating of 1500"

# Column width adjustments (closest achievable values given engine's pixel rounding)
$ws.Columns.Item(2).ColumnWidth = 16.42
$ws.Columns.Item(4).ColumnWidth = 73.92

# View state: zoom + final selection, matching the saved window state
$excel.ActiveWindow.Zoom = 115
$ws.Range("D22").Select()
